$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.682.71'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '1.646.35'
$ws.Range("E3").Value = '  -0.49%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '213.14'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '0.532'
$ws.Range("E6").Value = '  +3.34%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '23.06'
$ws.Range("E8").Value = '  -2.39%  '
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("E10").Value = '  -0.35%  '
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").Value = '1.878.49'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("D13").Value = '1.641.97'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").Value = '0.563'
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("D16").Value = '64.16'
$ws.Range("E16").Value = '  -2.42%  '
$ws.Range("D17").Value = '27.639.54'
$ws.Range("E17").Value = '  +1.15%  '
$ws.Range("D18").Value = '229.80'
$ws.Range("E18").Value = '  -1.02%  '
$ws.Range("D20").Value = '7.64'
$ws.Range("E20").Value = '  +2.51%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").Value = '4.32'
$ws.Range("E22").Value = '  -1.18%  '
$ws.Range("D23").Value = '10.03'
$ws.Range("E23").Value = '  +8.99%  '
$ws.Range("E24").Value = '  -2.81%  '
$ws.Range("D25").Value = '149.06'
$ws.Range("E25").Value = '  +1.12%  '
$ws.Range("D26").Value = '6.96'
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("E28").Value = '  +0.07%  '
$ws.Range("D29").Value = '15.64'
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = '1.19'
$ws.Range("E31").Value = '  -2.54%  '
$ws.Range("D32").Value = '3.30'
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +1.75%  '
$ws.Range("D34").Value = '1.441.09'
$ws.Range("E34").Value = '  -0.40%  '
$ws.Range("D35").Value = '1.59'
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").Value = '0.573'
$ws.Range("E37").Value = '  +0.31%  '
$ws.Range("D38").Value = '0.883'
$ws.Range("E38").Value = '  -2.92%  '
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").Value = '0.897'
$ws.Range("E40").Value = '  +14.36%  '
$ws.Range("E41").Value = '  -1.81%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +3.25%  '
$ws.Range("D44").Value = '2.47'
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("D46").Value = '65.51'
$ws.Range("E46").Value = '  +0.72%  '
$ws.Range("D47").Value = '1.787.59'
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("D48").Value = '1.69'
$ws.Range("E48").Value = '  -1.35%  '
$ws.Range("D49").Value = '86.44'
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").Value = '0.0989'
$ws.Range("E50").Value = '  -2.27%  '
$ws.Range("D51").Value = '7.78'
$ws.Range("E51").Value = '  +0.78%  '
